$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '42.297.22'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -1.06%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.249.64'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  +0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '247.64'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('E6').Value = '  -3.47%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '74.35'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -1.09%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -3.82%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '42.21'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +6.81%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0938'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -4.16%  '
$ws.Range('E12').Value = '  -3.79%  '
$ws.Range('E13').Value = '  -3.90%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '2.583.31'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.34%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '14.53'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -3.07%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.853'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -1.50%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '2.248.43'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.82%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '42.153.78'
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.0₃0979'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('E20').Value = '  -1.50%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '71.96'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('E22').Value = '  +4.54%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '230.25'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -2.61%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '8.26'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +30.18%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '11.10'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.55%  '
$ws.Range('E27').Value = '  -7.21%  '
$ws.Range('E28').Value = '  -3.55%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.16'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.28%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '169.31'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('E31').Value = '  -1.71%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.0822'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -5.95%  '
$ws.Range('E33').Value = '  -5.46%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '30.50'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -3.93%  '
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('E36').Value = '  +1.05%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '4.99'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +4.78%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.0308'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.87%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '13.47'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('E40').Value = '  -4.41%  '
$ws.Range('E41').Value = '  -1.50%  '
$ws.Range('E42').Value = '  -2.73%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '61.34'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.07%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '107.45'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('E45').Value = '  -3.36%  '
$ws.Range('E46').Value = '  +0.14%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.35%  '
$ws.Range('E48').Value = '  -3.53%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.16'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('E50').Value = '  +1.20%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '4.11'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -2.61%  '
